$wb = $excel.ActiveWorkbook


$ws = $wb.Worksheets.Item("ALC")
$ws.Range("H9").Value = 137.81818
$ws.Range("I9").Value = 133.25
$ws.Range("J9").Value = 150
$ws.Range("K9").Value = 133.25
$ws.Range("L9").Value = 150
$ws.Range("M9").Value = 35.75
$ws.Range("N9").Value = -488
$ws.Range("H18").Value = 12998.111
$ws.Range("J18").Value = 16997.334
$ws.Range("L18").Value = 16997.334
$ws.Range("N18").Value = -17565.334
$ws.Range("H40").Value = 2797
$ws.Range("I40").Value = 2666.6667
$ws.Range("J40").Value = 2927.3333
$ws.Range("K40").Value = 2666.6667
$ws.Range("L40").Value = 2927.3333
$ws.Range("M40").Value = -2491.6667
$ws.Range("N40").Value = -3277.3333
$ws.Range("H41").Value = 607.7646999999999
$ws.Range("I41").Value = 265
$ws.Range("J41").Value = 713.2308
$ws.Range("K41").Value = 265
$ws.Range("L41").Value = 713.2308
$ws.Range("M41").Value = 175
$ws.Range("N41").Value = -1593.2308
$ws.Range("H86").Value = 1187.5
$ws.Range("J86").Value = 1300
$ws.Range("L86").Value = 1300
$ws.Range("N86").Value = -3546
$ws.Range("H89").Value = 1187.5
$ws.Range("J89").Value = 1300
$ws.Range("L89").Value = 6500
$ws.Range("N89").Value = -17732
$ws.Range("H125").Value = 1244.6666
$ws.Range("I125").Value = 1214.3334
$ws.Range("J125").Value = 1335.6666
$ws.Range("K125").Value = 10929.0006
$ws.Range("L125").Value = 12020.9994
$ws.Range("M125").Value = -8469.000599999999
$ws.Range("N125").Value = -16940.9994
$ws.Range("H129").Value = 864.8421
$ws.Range("I129").Value = 649.3333
$ws.Range("J129").Value = 890.1961
$ws.Range("K129").Value = 1947.9999
$ws.Range("L129").Value = 2670.5883
$ws.Range("M129").Value = 3052.0001
$ws.Range("N129").Value = -12670.5883
$ws.Range("H137").Value = 1718.25
$ws.Range("I137").Value = 1363.3334
$ws.Range("J137").Value = 2250.625
$ws.Range("K137").Value = 4090.0002
$ws.Range("L137").Value = 6751.875
$ws.Range("M137").Value = -1540.0002
$ws.Range("N137").Value = -11851.875
$ws.Range("H138").Value = 1768.746
$ws.Range("I138").Value = 1274.4828
$ws.Range("J138").Value = 2190.3235
$ws.Range("K138").Value = 3823.4484
$ws.Range("L138").Value = 6570.970499999999
$ws.Range("M138").Value = 1316.5516
$ws.Range("N138").Value = -16850.9705

$ws = $wb.Worksheets.Item("ARM")
$ws.Range("H132").Value = 1604.4736
$ws.Range("I132").Value = 927.7143
$ws.Range("K132").Value = 2783.1429
$ws.Range("M132").Value = -253.1428999999998

$ws = $wb.Worksheets.Item("BSM")
$ws.Range("H86").Value = 102911.55
$ws.Range("I86").Value = 4372
$ws.Range("J86").Value = 168604.58
$ws.Range("K86").Value = 4372
$ws.Range("L86").Value = 168604.58
$ws.Range("M86").Value = -3249
$ws.Range("N86").Value = -170850.58
$ws.Range("H89").Value = 102911.55
$ws.Range("I89").Value = 4372
$ws.Range("J89").Value = 168604.58
$ws.Range("K89").Value = 21860
$ws.Range("L89").Value = 843022.8999999999
$ws.Range("M89").Value = -16244
$ws.Range("N89").Value = -854254.8999999999
$ws.Range("H108").Value = 94993.5
$ws.Range("J108").Value = 94993.5
$ws.Range("L108").Value = 94993.5
$ws.Range("N108").Value = -102673.5
$ws.Range("H134").Value = 2457.0715
$ws.Range("I134").Value = 2259.3438
$ws.Range("K134").Value = 6778.0314
$ws.Range("M134").Value = -4243.0314

$ws = $wb.Worksheets.Item("CRP")
$ws.Range("H22").Value = 1046
$ws.Range("I22").Value = 400
$ws.Range("J22").Value = 1207.5
$ws.Range("K22").Value = 400
$ws.Range("L22").Value = 1207.5
$ws.Range("M22").Value = -50
$ws.Range("N22").Value = -1907.5
$ws.Range("H31").Value = 2778.7
$ws.Range("I31").Value = 3160
$ws.Range("J31").Value = 2651.6
$ws.Range("K31").Value = 3160
$ws.Range("L31").Value = 2651.6
$ws.Range("M31").Value = -2865
$ws.Range("N31").Value = -3241.6
$ws.Range("H34").Value = 2778.7
$ws.Range("I34").Value = 3160
$ws.Range("J34").Value = 2651.6
$ws.Range("K34").Value = 3160
$ws.Range("L34").Value = 2651.6
$ws.Range("M34").Value = -2958
$ws.Range("N34").Value = -3055.6
$ws.Range("H50").Value = 14186.667
$ws.Range("J50").Value = 14186.667
$ws.Range("L50").Value = 14186.667
$ws.Range("N50").Value = -15436.667
$ws.Range("H59").Value = 21483.334
$ws.Range("J59").Value = 21483.334
$ws.Range("L59").Value = 21483.334
$ws.Range("N59").Value = -23773.334
$ws.Range("H60").Value = 16833.777
$ws.Range("J60").Value = 16833.777
$ws.Range("L60").Value = 16833.777
$ws.Range("N60").Value = -17855.777
$ws.Range("H62").Value = 2538.4
$ws.Range("I62").Value = 2538.4
$ws.Range("K62").Value = 2538.4
$ws.Range("M62").Value = -1914.4
$ws.Range("H65").Value = 2538.4
$ws.Range("I65").Value = 2538.4
$ws.Range("K65").Value = 12692
$ws.Range("M65").Value = -9572

$ws = $wb.Worksheets.Item("CUL")
$ws.Range("H17").Value = 4139.8
$ws.Range("I17").Value = 174.25
$ws.Range("J17").Value = 20002
$ws.Range("K17").Value = 522.75
$ws.Range("L17").Value = 60006
$ws.Range("M17").Value = -353.75
$ws.Range("N17").Value = -60344
$ws.Range("H32").Value = 1000
$ws.Range("I32").Value = 0
$ws.Range("K32").Value = 0
$ws.Range("M32").ClearContents()
$ws.Range("H56").Value = 7179.5
$ws.Range("I56").Value = 7179.5
$ws.Range("K56").Value = 7179.5
$ws.Range("M56").Value = -6649.5
$ws.Range("H61").Value = 320.625
$ws.Range("I61").Value = 575
$ws.Range("J61").Value = 235.83333
$ws.Range("K61").Value = 1725
$ws.Range("L61").Value = 707.49999
$ws.Range("M61").Value = -1510
$ws.Range("N61").Value = -1137.49999
$ws.Range("H130").Value = 3255.4443
$ws.Range("J130").Value = 3500
$ws.Range("L130").Value = 10500
$ws.Range("N130").Value = -20540
$ws.Range("H131").Value = 795.87
$ws.Range("J131").Value = 809.8936
$ws.Range("L131").Value = 2429.6808
$ws.Range("N131").Value = -12509.6808

$ws = $wb.Worksheets.Item("GSM")
$ws.Range("H97").Value = 1714.5714
$ws.Range("I97").Value = 1667
$ws.Range("J97").Value = 2000
$ws.Range("K97").Value = 1667
$ws.Range("L97").Value = 2000
$ws.Range("M97").Value = -1171
$ws.Range("N97").Value = -2992
$ws.Range("H122").Value = 2299.5
$ws.Range("I122").Value = 1779.4
$ws.Range("J122").Value = 2671
$ws.Range("K122").Value = 5338.200000000001
$ws.Range("L122").Value = 8013
$ws.Range("M122").Value = -2888.200000000001
$ws.Range("N122").Value = -12913
$ws.Range("H132").Value = 6413074.5
$ws.Range("I132").Value = 12822183
$ws.Range("J132").Value = 3966.3333
$ws.Range("K132").Value = 38466549
$ws.Range("L132").Value = 11898.9999
$ws.Range("M132").Value = -38464019
$ws.Range("N132").Value = -16958.9999

$ws = $wb.Worksheets.Item("LTW")
$ws.Range("H7").Value = 3404.5386
$ws.Range("I7").Value = 2205.818
$ws.Range("K7").Value = 2205.818
$ws.Range("M7").Value = -2093.818
$ws.Range("H20").Value = 9484.714
$ws.Range("I20").Value = 7278.6
$ws.Range("K20").Value = 7278.6
$ws.Range("M20").Value = -7052.6
$ws.Range("H22").Value = 3361.111
$ws.Range("I22").Value = 3950
$ws.Range("J22").Value = 2625
$ws.Range("K22").Value = 3950
$ws.Range("L22").Value = 2625
$ws.Range("M22").Value = -3655
$ws.Range("N22").Value = -3215
$ws.Range("H27").Value = 3361.111
$ws.Range("I27").Value = 3950
$ws.Range("J27").Value = 2625
$ws.Range("K27").Value = 3950
$ws.Range("L27").Value = 2625
$ws.Range("M27").Value = -3843
$ws.Range("N27").Value = -2839
$ws.Range("H68").Value = 3060
$ws.Range("I68").Value = 3106.8333
$ws.Range("J68").Value = 2966.3333
$ws.Range("K68").Value = 3106.8333
$ws.Range("L68").Value = 2966.3333
$ws.Range("M68").Value = -2357.8333
$ws.Range("N68").Value = -4464.3333
$ws.Range("H71").Value = 3060
$ws.Range("I71").Value = 3106.8333
$ws.Range("J71").Value = 2966.3333
$ws.Range("K71").Value = 15534.1665
$ws.Range("L71").Value = 14831.6665
$ws.Range("M71").Value = -11790.1665
$ws.Range("N71").Value = -22319.6665
$ws.Range("H126").Value = 3404.5386
$ws.Range("I126").Value = 2205.818
$ws.Range("K126").Value = 6617.454000000001
$ws.Range("M126").Value = -4147.454000000001
$ws.Range("H136").Value = 3933.8096
$ws.Range("I136").Value = 3043.5715
$ws.Range("K136").Value = 9130.7145
$ws.Range("M136").Value = -6580.7145

$ws = $wb.Worksheets.Item("WVR")
$ws.Range("H14").Value = 1057
$ws.Range("J14").Value = 1057
$ws.Range("L14").Value = 1057
$ws.Range("N14").Value = -1393
$ws.Range("H96").Value = 2320
$ws.Range("I96").Value = 0
$ws.Range("J96").Value = 2320
$ws.Range("K96").Value = 0
$ws.Range("L96").Value = 2320
$ws.Range("M96").ClearContents()
$ws.Range("N96").Value = -5066
$ws.Range("H107").Value = 820.9231
$ws.Range("I107").Value = 610.4286
$ws.Range("K107").Value = 1831.2858
$ws.Range("M107").Value = 88.71420000000012
